# Applies the cryptos.xlsx symbol-list refresh described in the commit:
# "Updated symbol list on Sat Jan 14 04:10:48 UTC 2023 with GitHub Actions"
# For every data row (2-51): Price (D) and Volume(1h) (E) are refreshed where
# the source changed, and Hora (G) moves from "3" to "4" for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) values, keyed by row number. Rows not listed keep their value.
$dNew = @{
    2 = "306.90"
    3 = "32.28"
    4 = "5.297"
    5 = "0.07444"
    6 = "7.759"
    7 = "3.714"
    8 = "1.563"
    9 = "0.9193"
    10 = "0.01632"
    11 = "0.1671"
    12 = "0.07526"
    13 = "0.07952"
    14 = "0.03101"
    15 = "0.09879"
    16 = "0.001528"
    17 = "0.04548"
    18 = "0.006134"
    19 = "3.472"
    20 = "2.243"
    21 = "0.3277"
    22 = "0.1325"
    23 = "4.238"
    25 = "0.001196"
    26 = "0.004543"
    40 = "0.04501"
    41 = "0.007287"
    42 = "0.1368"
    43 = "0.002257"
    44 = "0.01382"
    45 = "0.00006041"
    46 = "1.893"
}

# New Volume(1h) (column E) values, keyed by row number. Rows not listed keep their value.
$eNew = @{
    2 = "6.42%"
    3 = "9.69%"
    4 = "3.04%"
    5 = "11.61%"
    6 = "5.58%"
    7 = "9.16%"
    8 = "14.93%"
    9 = "0.19%"
    10 = "2,428.48%"
    11 = "5.92%"
    12 = "13.41%"
    13 = "4.41%"
    14 = "4.26%"
    15 = "9.77%"
    16 = "-4.11%"
    17 = "1.67%"
    18 = "-2.04%"
    19 = "0.39%"
    20 = "0.93%"
    21 = "1.97%"
    22 = "1.25%"
    23 = "4.23%"
    24 = "4.46%"
    25 = "0.55%"
    26 = "9.78%"
    27 = "-6.44%"
    28 = "2.91%"
    40 = "6.84%"
    41 = "7.92%"
    42 = "10.15%"
    43 = "14.09%"
    44 = "8.28%"
    45 = "7.89%"
    46 = "-3.80%"
    47 = "-0.57%"
}

function Set-TextValue($cell, $text) {
    # Force the cell to stay a plain text value (matches the source workbook,
    # which stores every cell as a literal/inline string) instead of letting
    # Excel auto-convert numeric- or percent-looking text into a number.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

for ($row = 2; $row -le 51; $row++) {
    if ($dNew.ContainsKey($row)) {
        Set-TextValue ($ws.Cells.Item($row, 4)) $dNew[$row]
    }
    if ($eNew.ContainsKey($row)) {
        Set-TextValue ($ws.Cells.Item($row, 5)) $eNew[$row]
    }
    # Hora (column G): every row goes from "3" to "4"
    Set-TextValue ($ws.Cells.Item($row, 7)) "4"
}

Write-Host "Updated Price/Volume/Hora for rows 2-51"
